$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest values pulled by the scheduled GitHub Actions job.
# A handful of Price cells look like plain numbers (e.g. "1.00", "8.29");
# a leading apostrophe keeps Excel from collapsing them to numeric values
# and losing the trailing zero / exact text formatting.

$ws.Range("D2").Value = "57.867.87"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.340.49"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'540.17"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'134.03"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  +5.41%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "'5.50"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'23.74"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "2.756.59"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "57.795.36"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "2.341.64"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "'329.02"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "'6.71"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'63.02"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'0.164"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'8.29"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("E27").Value = "  -5.39%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'170.05"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "'18.32"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "'4.17"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'39.08"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'141.38"
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'287.31"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "'0.0944"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "'19.14"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "'0.565"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'0.379"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'11.08"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'0.953"
$ws.Range("E51").Value = "  +0.94%  "
